$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update each changed cell with its new value from the refreshed crypto feed.
# A handful of Price cells in column D are plain decimal numbers (e.g. "300.09");
# the source data keeps them as literal text, so force a text number format
# before writing them (otherwise Excel auto-converts the string to a numeric
# value and introduces floating point rounding), then clear that temporary
# format again so the cell ends up with no explicit style, same as before.

$ws.Range('D2').Value = '42.929.44'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '2.303.08'
$ws.Range('E3').Value = '  -0.90%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.09'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.74'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.98%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -1.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.86'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.89%  '
$ws.Range('E12').Value = '  +0.66%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '17.71'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.79'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.23%  '
$ws.Range('D15').Value = '2.659.74'
$ws.Range('E15').Value = '  -0.91%  '
$ws.Range('D16').Value = '2.300.97'
$ws.Range('E16').Value = '  -0.67%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.779'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.47%  '
$ws.Range('D18').Value = '42.894.54'
$ws.Range('E18').Value = '  -0.46%  '
$ws.Range('E19').Value = '  -2.82%  '
$ws.Range('D20').Value = '0.0₃0907'
$ws.Range('E20').Value = '  -0.72%  '
$ws.Range('E21').Value = '  -2.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.97'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '242.03'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.48%  '
$ws.Range('E24').Value = '  -1.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.43'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.18%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '25.08'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '166.28'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.09%  '
$ws.Range('E30').Value = '  -0.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.04'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.82%  '
$ws.Range('E32').Value = '  -4.47%  '
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.75'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -4.18%  '
$ws.Range('E35').Value = '  -3.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.49'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.38%  '
$ws.Range('E37').Value = '  -0.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0686'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.93%  '
$ws.Range('E39').Value = '  -2.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.77'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.76'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.92%  '
$ws.Range('E42').Value = '  +0.40%  '
$ws.Range('D43').Value = '2.000.43'
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0285'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.49%  '
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.15'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.02%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.14'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.22%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.24'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.12%  '
$ws.Range('E48').Value = '  -2.44%  '
$ws.Range('D49').Value = '2.524.60'
$ws.Range('E49').Value = '  -0.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '53.38'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.34'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -5.33%  '
